# Updated cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.979.23"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.742.51"
$ws.Range("E3").Value = "  -3.78%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.07"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3763"
$ws.Range("E7").Value = "  -4.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3347"
$ws.Range("E8").Value = "  -4.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.46"
$ws.Range("E9").Value = "  -6.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.112"
$ws.Range("E10").Value = "  -6.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07207"

$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.154"
$ws.Range("E14").Value = "  -5.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.101"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.748.50"
$ws.Range("E16").Value = "  -3.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001056"
$ws.Range("E17").Value = "  -4.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06571"
$ws.Range("E18").Value = "  -2.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.59"
$ws.Range("E19").Value = "  -6.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.81"
$ws.Range("E21").Value = "  -5.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.231"
$ws.Range("E22").Value = "  -5.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.986.20"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.63"
$ws.Range("E24").Value = "  -6.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.67"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.78"
$ws.Range("E27").Value = "  -7.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.320"
$ws.Range("E28").Value = "  -8.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.949.47"
$ws.Range("E29").Value = "  -3.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.87"
$ws.Range("E30").Value = "  -3.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.243"
$ws.Range("E31").Value = "  -16.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.775"
$ws.Range("E33").Value = "  -9.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08722"
$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.18"
$ws.Range("E35").Value = "  -7.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6676"
$ws.Range("E36").Value = "  -4.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02314"
$ws.Range("E37").Value = "  -6.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06195"
$ws.Range("E38").Value = "  -5.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.128"
$ws.Range("E39").Value = "  -6.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2106"
$ws.Range("E40").Value = "  -5.81%  "

$ws.Range("E41").Value = "  -4.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.439"
$ws.Range("E42").Value = "  -10.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.939"
$ws.Range("E44").Value = "  -7.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("E45").Value = "  -6.94%  "

$ws.Range("E46").Value = "  -1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6035"
$ws.Range("E47").Value = "  -6.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.77"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("E49").Value = "  -7.05%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.176"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07123"
$ws.Range("E51").Value = "  -1.65%  "
